$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1000
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 1000
$ws.Range("D2").Value = 1000
$ws.Range("G2").Value = 1000
